$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue 2 4 '51.787.94'
Set-TextValue 2 5 '  -1.06%  '
Set-TextValue 3 4 '2.930.39'
Set-TextValue 3 5 '  +0.47%  '
Set-TextValue 4 5 '  +0.09%  '
Set-TextValue 5 4 '353.24'
Set-TextValue 5 5 '  -0.36%  '
Set-TextValue 6 4 '107.88'
Set-TextValue 6 5 '  -5.94%  '
Set-TextValue 7 5 '  +0.87%  '
Set-TextValue 8 5 '  +0.05%  '
Set-TextValue 9 4 '0.617'
Set-TextValue 9 5 '  -1.10%  '
Set-TextValue 10 4 '37.98'
Set-TextValue 10 5 '  -5.01%  '
Set-TextValue 11 5 '  +1.01%  '
Set-TextValue 12 5 '  -1.03%  '
Set-TextValue 13 4 '19.14'
Set-TextValue 13 5 '  -3.80%  '
Set-TextValue 14 4 '3.394.93'
Set-TextValue 14 5 '  +0.80%  '
Set-TextValue 15 5 '  -0.66%  '
Set-TextValue 16 4 '2.945.34'
Set-TextValue 16 5 '  +2.22%  '
Set-TextValue 17 4 '0.964'
Set-TextValue 17 5 '  -2.59%  '
Set-TextValue 18 4 '51.762.90'
Set-TextValue 18 5 '  -1.15%  '
Set-TextValue 19 5 '  +2.26%  '
Set-TextValue 20 4 '7.48'
Set-TextValue 20 5 '  -2.15%  '
Set-TextValue 21 5 '  -4.35%  '
Set-TextValue 22 4 '0.0₃0969'
Set-TextValue 22 5 '  -1.17%  '
Set-TextValue 23 4 '69.54'
Set-TextValue 23 5 '  -2.41%  '
Set-TextValue 24 4 '265.04'
Set-TextValue 24 5 '  -1.91%  '
Set-TextValue 25 5 '  -2.82%  '
Set-TextValue 26 5 '  -3.57%  '
Set-TextValue 27 5 '  -0.39%  '
Set-TextValue 28 4 '7.55'
Set-TextValue 28 5 '  +11.60%  '
Set-TextValue 29 5 '  +0.18%  '
Set-TextValue 30 5 '  -0.12%  '
Set-TextValue 31 4 '10.26'
Set-TextValue 31 5 '  -4.01%  '
Set-TextValue 32 4 '36.50'
Set-TextValue 32 5 '  -3.51%  '
Set-TextValue 33 4 '2.17'
Set-TextValue 33 5 '  -4.90%  '
Set-TextValue 34 5 '  -2.71%  '
Set-TextValue 35 4 '51.88'
Set-TextValue 35 5 '  -2.40%  '
Set-TextValue 36 4 '0.0429'
Set-TextValue 36 5 '  -4.81%  '
Set-TextValue 37 5 '  +0.05%  '
Set-TextValue 38 4 '3.14'
Set-TextValue 38 5 '  -6.21%  '
Set-TextValue 39 5 '  -5.78%  '
Set-TextValue 40 5 '  -4.11%  '
Set-TextValue 41 4 '2.66'
Set-TextValue 41 5 '  -3.18%  '
Set-TextValue 42 5 '  -0.24%  '
Set-TextValue 43 4 '22.95'
Set-TextValue 43 5 '  -0.86%  '
Set-TextValue 44 4 '118.67'
Set-TextValue 44 5 '  -0.26%  '
Set-TextValue 45 5 '  -0.90%  '
Set-TextValue 46 5 '  -2.98%  '
Set-TextValue 47 4 '2.112.70'
Set-TextValue 47 5 '  -3.32%  '
Set-TextValue 48 5 '  -5.58%  '
Set-TextValue 49 2 'TheGraph'
Set-TextValue 49 3 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue 49 4 '0.239'
Set-TextValue 49 5 '  -8.84%  '
Set-TextValue 50 2 'BEAM'
Set-TextValue 50 3 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
Set-TextValue 50 4 '0.0339'
Set-TextValue 50 5 '  -3.68%  '
Set-TextValue 51 2 'SEI'
Set-TextValue 51 3 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
Set-TextValue 51 4 '0.903'
Set-TextValue 51 5 '  -5.99%  '
